# Add new notebooks for experimental datasets:
#  - Populate a new "Non-Linear fit success" (Yes/No/NA) column of data in
#    column E (the header text itself already existed).
#  - Fix the spelling of the "Significant Order (linear epistasis)" header
#    in column D.
#  - Move the active selection to F2.
#  - Update the workbook's recorded absolute folder path.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-affirm the E1 header text (unchanged wording) before introducing the
# brand new "Yes"/"No"/"NA" strings, then fix the misspelled D1 header, so
# that the shared-string table ends up ordered the same way the authored
# workbook has it.
$ws.Range("E1").Value = "Non-Linear fit success"

$ws.Range("E2").Value = "Yes"

$ws.Range("D1").Value = "Significant Order (linear epistasis)"

$ws.Range("E3").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("E5").Value = "NA"

$ws.Range("E6").Value = "No"
$ws.Range("E7").Value = "No"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"

$ws.Range("E10").Value = "Yes"
$ws.Range("E11").Value = "Yes"
$ws.Range("E12").Value = "Yes"
$ws.Range("E13").Value = "Yes"

# Move the current selection to F2.
$ws.Range("F2").Select()

# Update the workbook's recorded absolute folder path (Mac Excel metadata).
$wb.AbsolutePath = "/Users/Zsailer/Documents/Research/projects/notebooks/epistasis-notebooks/"
